$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1155809658704365
$ws.Range("C2").Value = 0.390457646890028
$ws.Range("D2").Value = 0.3755683084142444
$ws.Range("E2").Value = 0.612836281901002
$ws.Range("F2").Value = 0.6153644530269748

$ws.Range("B3").Value = 0.6511545689434494
$ws.Range("C3").Value = 0.8016979307476643
$ws.Range("D3").Value = 4.674211074126436
$ws.Range("E3").Value = 2.161992385307228
$ws.Range("F3").Value = 2.107937406714976
$ws.Range("G3").Value = 23

$ws.Range("B4").Value = 0.2201423283843545
$ws.Range("C4").Value = 1.112030744985122
$ws.Range("D4").Value = 7.605937587058003
$ws.Range("E4").Value = 2.757886434764492
$ws.Range("F4").Value = 2.810871133758786
$ws.Range("G4").Value = 23

$ws.Range("B5").Value = 0.2501754929190889
$ws.Range("C5").Value = 1.125383142428414
$ws.Range("D5").Value = 7.653888439551237
$ws.Range("E5").Value = 2.766566182029853
$ws.Range("F5").Value = 2.817154554387687
$ws.Range("G5").Value = 23

$ws.Range("B6").Value = 0.3004374354251687
$ws.Range("C6").Value = 1.183046795173429
$ws.Range("D6").Value = 7.775752270546474
$ws.Range("E6").Value = 2.788503589839266
$ws.Range("F6").Value = 2.834577579661701
$ws.Range("G6").Value = 23

$ws.Range("B7").Value = 0.2656501924638524
$ws.Range("C7").Value = 1.272698581772992
$ws.Range("D7").Value = 7.926928235987235
$ws.Range("E7").Value = 2.815480107546
$ws.Range("F7").Value = 2.86591440951239
$ws.Range("G7").Value = 23

$ws.Range("B8").Value = 0.1735487773573387
$ws.Range("C8").Value = 1.27905897888513
$ws.Range("D8").Value = 7.992029374264993
$ws.Range("E8").Value = 2.827017752732549
$ws.Range("F8").Value = 2.885102286758365
$ws.Range("G8").Value = 23

$ws.Range("B9").Value = 0.2219982892021358
$ws.Range("C9").Value = 1.341063655789804
$ws.Range("D9").Value = 8.0376137195044
$ws.Range("E9").Value = 2.835068556402896
$ws.Range("F9").Value = 2.889885189754726
$ws.Range("G9").Value = 23

$ws.Range("B10").Value = 0.1895698287473124
$ws.Range("C10").Value = 1.293626197483444
$ws.Range("D10").Value = 8.015968031220561
$ws.Range("E10").Value = 2.831248493371884
$ws.Range("F10").Value = 2.88838363227864
$ws.Range("G10").Value = 23

$ws.Range("B11").Value = 0.08322311406686993
$ws.Range("C11").Value = 1.20051888709089
$ws.Range("D11").Value = 7.697829717014696
$ws.Range("E11").Value = 2.77449629969382
$ws.Range("F11").Value = 2.835575807300819
$ws.Range("G11").Value = 23
